# Refresh the cryptos price list (Price / Volume(1h) columns) with the
# latest scraped values. Numeric-looking "Price" strings are written with a
# leading apostrophe so Excel keeps them as text (matching the original
# inline-string cell type) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.446.01'
$ws.Range('E2').Value = '  -1.37%  '
$ws.Range('D3').Value = '1.686.38'
$ws.Range('D4').Value = '''0.9994'
$ws.Range('E4').Value = '  -0.44%  '
$ws.Range('D5').Value = '''316.16'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '''0.9988'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').Value = '''0.3896'
$ws.Range('E7').Value = '  -1.04%  '
$ws.Range('D8').Value = '''0.4025'
$ws.Range('E8').Value = '  -0.41%  '
$ws.Range('D9').Value = '''1.486'
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').Value = '''0.9995'
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('D11').Value = '''52.43'
$ws.Range('E11').Value = '  -2.62%  '
$ws.Range('D12').Value = '''0.08763'
$ws.Range('E12').Value = '  -1.55%  '
$ws.Range('D13').Value = '''25.96'
$ws.Range('E13').Value = '  +10.72%  '
$ws.Range('D14').Value = '''7.482'
$ws.Range('E14').Value = '  +3.24%  '
$ws.Range('D15').Value = '''8.124'
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = '''0.00001348'
$ws.Range('E16').Value = '  +1.41%  '
$ws.Range('D17').Value = '1.684.45'
$ws.Range('E17').Value = '  -1.21%  '
$ws.Range('D18').Value = '''97.96'
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('D19').Value = '''0.07264'
$ws.Range('E19').Value = '  +3.13%  '
$ws.Range('D20').Value = '''19.88'
$ws.Range('E20').Value = '  +1.14%  '
$ws.Range('E21').Value = '  +3.28%  '
$ws.Range('D22').Value = '''0.9984'
$ws.Range('E22').Value = '  -0.41%  '
$ws.Range('E23').Value = '  -2.50%  '
$ws.Range('D24').Value = '24.429.64'
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('D25').Value = '''3.027'
$ws.Range('E25').Value = '  -5.50%  '
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('D27').Value = '''22.56'
$ws.Range('D28').Value = '''167.29'
$ws.Range('E28').Value = '  +3.20%  '
$ws.Range('D29').Value = '''8.637'
$ws.Range('E29').Value = '  +7.98%  '
$ws.Range('D30').Value = '''5.353'
$ws.Range('E30').Value = '  +3.45%  '
$ws.Range('D31').Value = '''138.42'
$ws.Range('E31').Value = '  +1.39%  '
$ws.Range('D32').Value = '1.868.67'
$ws.Range('E32').Value = '  -1.20%  '
$ws.Range('D33').Value = '''0.08757'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('D34').Value = '''7.336'
$ws.Range('E34').Value = '  +0.57%  '
$ws.Range('B35').Value = 'WEMIXTOKEN'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '''2.110'
$ws.Range('E35').Value = '  +6.72%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '''1.047'
$ws.Range('E36').Value = '  -2.89%  '
$ws.Range('D37').Value = '''0.03019'
$ws.Range('E37').Value = '  +9.69%  '
$ws.Range('D38').Value = '''0.2782'
$ws.Range('D39').Value = '''10.81'
$ws.Range('E39').Value = '  -3.29%  '
$ws.Range('D40').Value = '''0.09132'
$ws.Range('E40').Value = '  -0.80%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').Value = '''0.8040'
$ws.Range('E41').Value = '  +4.42%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '''14.14'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('D43').Value = '''1.472'
$ws.Range('E43').Value = '  +0.70%  '
$ws.Range('D44').Value = '''17.51'
$ws.Range('E44').Value = '  +10.75%  '
$ws.Range('D45').Value = '''2.625'
$ws.Range('E45').Value = '  +2.05%  '
$ws.Range('D46').Value = '''0.7236'
$ws.Range('E46').Value = '  +0.90%  '
$ws.Range('D47').Value = '''4.264'
$ws.Range('E47').Value = '  +1.36%  '
$ws.Range('D48').Value = '''1.426'
$ws.Range('E48').Value = '  +8.66%  '
$ws.Range('D49').Value = '''0.9987'
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('D50').Value = '''139.06'
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('D51').Value = '''0.08074'
$ws.Range('E51').Value = '  +0.97%  '
